$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column I (rows 4-14) into the new column J, then set the 2021 values.
$rows = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14)
$values = @{
    4  = 2021
    5  = 24.4
    6  = 45.7
    7  = 38
    8  = 51.3
    9  = 51.5
    10 = 13
    11 = 36.4
    12 = 27
    13 = 2.7
    14 = 40.4
}

foreach ($r in $rows) {
    $src = $ws.Cells.Item($r, 9)   # column I
    $dst = $ws.Cells.Item($r, 10)  # column J
    [void]$src.Copy()
    [void]$dst.PasteSpecial(-4122) # xlPasteFormats
    $dst.Value = $values[$r]
}

# Row 3 height changes as part of this edit.
$ws.Rows.Item(3).RowHeight = 13.5

# Update the selected cell to match the saved view state.
[void]$ws.Range("K18").Select()
